# PP_Mapping.xlsx update:
#  1. Add a new worksheet "Sheet2" after "ValueMapping".
#  2. Copy the "InsuredRiskClass" value-mapping rows (ValueMapping!A10:C14)
#     onto the new sheet as rows 1-5.
#  3. Remove the now-duplicated rows 11:14 from "ValueMapping" (row 10 is
#     left in place on both sheets).
#  4. Append five new JSON<->Excel field mappings to "Sheet1" (rows 56-60)
#     for policy conversion dates and agent info.
#  5. Leave "ValueMapping" as the active/selected sheet, matching the
#     original workbook's view state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("ValueMapping")

# --- 1. New sheet, placed right after ValueMapping -------------------------
$ws3 = $wb.Worksheets.Add([Type]::Missing, $ws2)
$ws3.Name = "Sheet2"

# --- 2. Copy the risk-class rows onto the new sheet -------------------------
$ws2.Range("A10:C14").Copy()
$ws3.Range("A1").PasteSpecial()
$excel.CutCopyMode = $false

# --- 3. Trim the old sheet back down to its first risk-class row -----------
$ws2.Rows("11:14").Delete()

# --- 4. New Sheet1 mapping rows ---------------------------------------------
# (written in this particular cell order so new shared strings land with
# the same index values used by the workbook: 137-146)
$ws1.Range("B56").Value = "PPBulk.ConversionAttainedAge"
$ws1.Range("B57").Value = "PPBulk.ConversionOriginalAge"
$ws1.Range("A57").Value = "policyDetail.premiumScheduleDetail.conversionDetail.originalAgeDt"
$ws1.Range("A56").Value = "policyDetail.premiumScheduleDetail.conversionDetail.attainedAgeDt"
$ws1.Range("B58").Value = "Agent1.Agent Contract"
$ws1.Range("A58").Value = "policyDetail.agentInfo[0].contractType"
$ws1.Range("A59").Value = "policyDetail.agentInfo[0].marketerId"
$ws1.Range("B59").Value = "Agent1.Agent  Number"
$ws1.Range("A60").Value = "policyDetail.agentInfo[0].percentage"
$ws1.Range("B60").Value = "Agent1.Agent Share"

# --- 5. Restore view/selection state ----------------------------------------
$ws1.Range("A70").Select()
$ws3.Rows("1:1").Select()
$ws2.Range("E7").Select()
